# Oppdaterer kandidatproduksjonstall (Tilbud/Differanse) for Lektorutdannede,
# Faglærere og Yrkesfaglærere (rader 45-64, 66-85, 87-106) med nye
# innfilsverdier for Kandidatproduksjon.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C45").Value = 5963
$ws.Range("E45").Value = 421
$ws.Range("C46").Value = 6431
$ws.Range("E46").Value = 870
$ws.Range("C47").Value = 6896
$ws.Range("E47").Value = 1322
$ws.Range("C48").Value = 7360
$ws.Range("E48").Value = 1777
$ws.Range("C49").Value = 7821
$ws.Range("E49").Value = 2227
$ws.Range("C50").Value = 8272
$ws.Range("E50").Value = 2672
$ws.Range("C51").Value = 8718
$ws.Range("E51").Value = 3118
$ws.Range("C52").Value = 9151
$ws.Range("E52").Value = 3569
$ws.Range("C53").Value = 9578
$ws.Range("E53").Value = 4021
$ws.Range("C54").Value = 9990
$ws.Range("E54").Value = 4459
$ws.Range("C55").Value = 10396
$ws.Range("E55").Value = 4884
$ws.Range("C56").Value = 10790
$ws.Range("E56").Value = 5291
$ws.Range("C57").Value = 11177
$ws.Range("E57").Value = 5684
$ws.Range("C58").Value = 11561
$ws.Range("E58").Value = 6081
$ws.Range("C59").Value = 11940
$ws.Range("E59").Value = 6477
$ws.Range("C60").Value = 12308
$ws.Range("E60").Value = 6859
$ws.Range("C61").Value = 12677
$ws.Range("E61").Value = 7225
$ws.Range("C62").Value = 13044
$ws.Range("E62").Value = 7583
$ws.Range("C63").Value = 13409
$ws.Range("E63").Value = 7936
$ws.Range("C64").Value = 13769
$ws.Range("E64").Value = 8280
$ws.Range("C66").Value = 5812
$ws.Range("E66").Value = 271
$ws.Range("C67").Value = 6122
$ws.Range("E67").Value = 561
$ws.Range("C68").Value = 6427
$ws.Range("E68").Value = 852
$ws.Range("C69").Value = 6728
$ws.Range("E69").Value = 1145
$ws.Range("C70").Value = 7026
$ws.Range("E70").Value = 1433
$ws.Range("C71").Value = 7315
$ws.Range("E71").Value = 1715
$ws.Range("C72").Value = 7600
$ws.Range("E72").Value = 2001
$ws.Range("C73").Value = 7874
$ws.Range("E73").Value = 2293
$ws.Range("C74").Value = 8144
$ws.Range("E74").Value = 2587
$ws.Range("C75").Value = 8401
$ws.Range("E75").Value = 2870
$ws.Range("C76").Value = 8652
$ws.Range("E76").Value = 3140
$ws.Range("C77").Value = 8894
$ws.Range("E77").Value = 3395
$ws.Range("C78").Value = 9130
$ws.Range("E78").Value = 3636
$ws.Range("C79").Value = 9363
$ws.Range("E79").Value = 3883
$ws.Range("C80").Value = 9594
$ws.Range("E80").Value = 4131
$ws.Range("C81").Value = 9816
$ws.Range("E81").Value = 4366
$ws.Range("C82").Value = 10039
$ws.Range("E82").Value = 4588
$ws.Range("C83").Value = 10264
$ws.Range("E83").Value = 4803
$ws.Range("C84").Value = 10487
$ws.Range("E84").Value = 5014
$ws.Range("C85").Value = 10707
$ws.Range("E85").Value = 5218
$ws.Range("C87").Value = 5633
$ws.Range("E87").Value = 92
$ws.Range("C88").Value = 5781
$ws.Range("E88").Value = 220
$ws.Range("C89").Value = 5927
$ws.Range("E89").Value = 352
$ws.Range("C90").Value = 6078
$ws.Range("E90").Value = 495
$ws.Range("C91").Value = 6231
$ws.Range("E91").Value = 637
$ws.Range("C92").Value = 6378
$ws.Range("E92").Value = 778
$ws.Range("C93").Value = 6524
$ws.Range("E93").Value = 924
$ws.Range("C94").Value = 6662
$ws.Range("E94").Value = 1080
$ws.Range("C95").Value = 6798
$ws.Range("E95").Value = 1241
$ws.Range("C96").Value = 6924
$ws.Range("E96").Value = 1393
$ws.Range("C97").Value = 7045
$ws.Range("E97").Value = 1533
$ws.Range("C98").Value = 7156
$ws.Range("E98").Value = 1657
$ws.Range("C99").Value = 7263
$ws.Range("E99").Value = 1769
$ws.Range("C100").Value = 7368
$ws.Range("E100").Value = 1887
$ws.Range("C101").Value = 7472
$ws.Range("E101").Value = 2009
$ws.Range("C102").Value = 7572
$ws.Range("E102").Value = 2122
$ws.Range("C103").Value = 7674
$ws.Range("E103").Value = 2223
$ws.Range("C104").Value = 7777
$ws.Range("E104").Value = 2315
$ws.Range("C105").Value = 7881
$ws.Range("E105").Value = 2407
$ws.Range("C106").Value = 7983
$ws.Range("E106").Value = 2494

